# Updates on Covid-19 cases in Kenya as of 30th April 2020
# Adds a new data row (row 48) for 30-Apr-2020 to the Sheet1 dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the new row's values ---
# Columns: A Date, B New Cases, C Tested, D Travelled From, E County,
#          F Aggregation, G Case Type, H Recover, I Death,
#          J URL, K Info Giver, L Ages, M Age Recover, N Age Death,
#          O Female Infected, P Male Infected
$ws.Range("A48").Value2 = 43951
$ws.Range("B48").Value2 = 12
$ws.Range("C48").Value2 = 777
$ws.Range("D48").Value2 = "None"
$ws.Range("E48").Value2 = "Mombasa(7),Nairobi(5)"
$ws.Range("F48").Value2 = 396
$ws.Range("G48").Value2 = "Community(12)"
$ws.Range("H48").Value2 = 15
$ws.Range("I48").Value2 = 2
# Set L before K so new shared strings are appended in the same order as the target workbook
$ws.Range("L48").Value2 = "1-75."
$ws.Range("K48").Value2 = "Aman"
$ws.Range("O48").Value2 = 3
$ws.Range("P48").Value2 = 9

# A48 should keep the date number format / alignment used by the rest of column A
$ws.Range("A48").NumberFormat = $ws.Range("A47").NumberFormat
$ws.Range("A48").HorizontalAlignment = $ws.Range("A47").HorizontalAlignment

# --- Update the sheet view to match where the author was working ---
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 35
[void]$ws.Range("K46").Select()
